$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add column K ("X"/present?) checkbox to the rows that gained an extra
# attendance mark (this raises their Q "total presences" sum by 1 via the
# existing shared SUM formula already present in column Q).
$rowsGainingK = @(7,8,9,12,14,15,16,17,18,25,26,29,37,38,42,47,49,54)
foreach ($r in $rowsGainingK) {
    $ws.Cells.Item($r, 11).Value = $true   # column K = 11
}

# Append a brand new student row (59) at the bottom of the table.
$ws.Range("B59").Value = "Cristea Octavian"

# Give B59 the same banded look as the other "first/odd" rows (fillId 3)
# but with only a left border, matching a new row added past the last
# grid line of the table: copy the light-fill style from row 4's name
# cell, then strip the top/bottom borders that B4 has.
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B59").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B59").Borders.Item(8).LineStyle = -4142  # xlEdgeTop, xlLineStyleNone
$ws.Range("B59").Borders.Item(9).LineStyle = -4142  # xlEdgeBottom, xlLineStyleNone
$ws.Application.CutCopyMode = $false

$ws.Cells.Item(59, 11).Value = $true   # K59 = TRUE

# Q59 total, reusing the same style as the other total cells in that
# region of the table (Q57/Q58 use style index 34).
$ws.Range("Q58").Copy() | Out-Null
$ws.Range("Q59").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("Q59").Formula = "=C59+D59+E59+F59+G59+H59+I59+J59+K59+L59+M59+N59+O59+P59"

# Extend the "below 10 presences" conditional formatting highlight down
# to the newly added row.
$fc = $ws.Range("Q3:Q58").FormatConditions
$cond = $fc.Item(1)
$cond.ModifyAppliesToRange($ws.Range("Q3:Q59"))

# Keep the bottom-right frozen pane's active selection pointed at the new
# row, mirroring where the user ended up after adding the entry.
$ws.Range("D61").Select()
